$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update price for the "塞纳里奥头" sale row (row 28)
$ws.Range("C28").Value = 100

# Correct the buyer name for that sale
$ws.Range("D28").Value = "马果果"

# Recalculate all formulas so dependent totals refresh
$excel.CalculateFullRebuild()

# Update the active selection to reflect where the edit was made
$ws.Range("D28").Select()
